# Add the new "Gigabit LAN Chip" BOM line (row 16) to the Bill of Materials
# sheet, mirroring the formatting already used by the other data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of BOM data directly below the last existing item (row 15).
$ws.Range("C16").Value2 = 10
$ws.Range("D16").Value2 = "Gigabit LAN Chip"
$ws.Range("E16").Value2 = "Gigabit Platform LAN Connect"
$ws.Range("F16").Value2 = "Intel"
$ws.Range("G16").Value2 = 82566

# Manufacturer P/N column in the other rows is left-aligned; match that here.
# -4131 == xlLeft
$ws.Range("G16").HorizontalAlignment = -4131

# Match the selection left behind by the edit.
$ws.Range("E17").Select() | Out-Null
